$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$subs = @(
    "sub_011", "sub_012", "sub_013", "sub_014", "sub_015",
    "sub_016", "sub_017", "sub_018", "sub_019", "sub_020"
)

$row = 12
foreach ($s in $subs) {
    $ws.Range("A$row").Value = $s
    $ws.Range("B$row").Value = $false
    $row = $row + 1
}

$ws.Range("B19").Select() | Out-Null
